$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the bold/bordered/centered header style (font 1, border 1, cellXf 1)
# so header cells revert to the default style, and clear A1's label text.
$ws.Range("A1:Z1").ClearFormats()
$ws.Range("A1").ClearContents()

# Corrected pre/post/total fixation metric values (rows 3-7)
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 19
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 17
$ws.Range("I3").Value = 18
$ws.Range("M3").Value = 34
$ws.Range("N3").Value = 17
$ws.Range("O3").Value = 3
$ws.Range("P3").Value = 16
$ws.Range("T3").Value = 33
$ws.Range("U3").Value = 6
$ws.Range("C4").Value = 2
$ws.Range("E4").Value = 131
$ws.Range("F4").Value = 36
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 44
$ws.Range("M4").Value = 101
$ws.Range("N4").Value = 22
$ws.Range("O4").Value = 5
$ws.Range("P4").Value = 24
$ws.Range("T4").Value = 75
$ws.Range("U4").Value = 7
$ws.Range("C5").Value = 817.5599999999999
$ws.Range("E5").Value = 47977.42
$ws.Range("F5").Value = 16234.07
$ws.Range("G5").Value = 6957.28
$ws.Range("H5").Value = 14148.51
$ws.Range("I5").Value = 16734.22
$ws.Range("M5").Value = 38390.9
$ws.Range("N5").Value = 10261.21
$ws.Range("O5").Value = 1985.37
$ws.Range("P5").Value = 10711.55
$ws.Range("T5").Value = 29031.5
$ws.Range("U5").Value = 3720.89
$ws.Range("B6").Value = 0.89
$ws.Range("C6").Value = 0.8
$ws.Range("D6").Value = 6.73
$ws.Range("E6").Value = 46.76
$ws.Range("F6").Value = 15.82
$ws.Range("G6").Value = 6.78
$ws.Range("H6").Value = 13.79
$ws.Range("I6").Value = 16.31
$ws.Range("J6").Value = 2.91
$ws.Range("K6").Value = 0.57
$ws.Range("M6").Value = 37.42
$ws.Range("N6").Value = 10
$ws.Range("O6").Value = 1.94
$ws.Range("P6").Value = 10.44
$ws.Range("Q6").Value = 3.82
$ws.Range("R6").Value = 11.51
$ws.Range("S6").Value = 0.33
$ws.Range("T6").Value = 28.3
$ws.Range("U6").Value = 3.63
$ws.Range("V6").Value = 0.24
$ws.Range("W6").Value = 2.75
$ws.Range("X6").Value = 0.23
$ws.Range("Y6").Value = 0.57
$ws.Range("Z6").Value = 2.18
$ws.Range("C7").Value = 408.78
$ws.Range("E7").Value = 366.24
$ws.Range("F7").Value = 450.95
$ws.Range("G7").Value = 409.25
$ws.Range("H7").Value = 393.01
$ws.Range("I7").Value = 380.32
$ws.Range("M7").Value = 380.11
$ws.Range("N7").Value = 466.42
$ws.Range("O7").Value = 397.07
$ws.Range("P7").Value = 446.31
$ws.Range("T7").Value = 387.09
$ws.Range("U7").Value = 531.5599999999999

# Drop the two trailing blank rows (10-11), shrinking the used range to A1:Z9
$ws.Rows("10:11").Delete()
